$wb = $excel.ActiveWorkbook
$episodios = $wb.Worksheets.Item("episodios")
$temporadas = $wb.Worksheets.Item("temporadas")

# ---------------------------------------------------------------------------
# temporadas: add "Temporada 2" as a new row (id_temporada, numero, titulo,
# descripcion, preview_img_url). Re-use the same season-1 preview image as
# the source workbook does.
# ---------------------------------------------------------------------------
$temporadas.Range("A3").Value = "Temproada_2"
$temporadas.Range("B3").Value = 2
$temporadas.Range("C3").Value = "Temporada 2"
$temporadas.Range("D3").Value = "Descripción Temporada 2"
$temporadas.Range("E3").Value = "Temporada_1.png"

# ---------------------------------------------------------------------------
# episodios: add a trailer row for season 2, cloned from the season-1
# trailer row (row 2) but pointing at "Temproada_2".
# ---------------------------------------------------------------------------
$episodios.Range("A8").Value = "trailer-1_temporada-1"
$episodios.Range("B8").Value = "Temproada_2"
$episodios.Range("C8").Value = "Trailer"
$episodios.Range("D8").Value = 0
$episodios.Range("E8").Value = "Trailer"
$episodios.Range("F8").Value = "https://collaboration.merck.com/sites/onconceptos/_layouts/15/embed.aspx?UniqueId=3934dfb1-e42b-4782-8a02-093fbd8326a6&embed=%7B%22hvm%22%3Atrue%2C%22ust%22%3Atrue%7D&referrer=OneUpFileViewer&referrerScenario=EmbedDialog.Create"
$episodios.Range("G8").Value = "https://images.pexels.com/photos/5726788/pexels-photo-5726788.jpeg?auto=compress&cs=tinysrgb&w=1920&h=1080&dpr=1"
$episodios.Range("H8").Value = 2023
$episodios.Range("I8").Value = "1h 35min"
$episodios.Range("J8").Value = "https://images.pexels.com/photos/5726788/pexels-photo-5726788.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=1"
$episodios.Range("K8").Value = "no"
$episodios.Range("L8").Value = "6/jun/2023"
$episodios.Range("M8").Value = 16
$episodios.Range("N8").Value = 20

# Give the new row the same look & feel as the row it was cloned from
# (general style + the hyperlink style used on column F).
$episodios.Range("A2:N2").Copy() | Out-Null
$episodios.Range("A8:N8").PasteSpecial(-4122) | Out-Null

# Re-create the hyperlink on F8 (same target/display as F2's).
$episodios.Hyperlinks.Add($episodios.Range("F8"), "https://collaboration.merck.com/sites/onconceptos/_layouts/15/embed.aspx?UniqueId=3934dfb1-e42b-4782-8a02-093fbd8326a6&embed=%7B%22hvm%22%3Atrue%2C%22ust%22%3Atrue%7D&referrer=OneUpFileViewer&referrerScenario=EmbedDialog.Create", "", "", "https://collaboration.merck.com/sites/onconceptos/_layouts/15/embed.aspx?UniqueId=3934dfb1-e42b-4782-8a02-093fbd8326a6&embed=%7B%22hvm%22%3Atrue%2C%22ust%22%3Atrue%7D&referrer=OneUpFileViewer&referrerScenario=EmbedDialog.Create") | Out-Null

# Restore F8's formatting/value after the hyperlink re-applied its own style.
$episodios.Range("F2").Copy() | Out-Null
$episodios.Range("F8").PasteSpecial(-4122) | Out-Null
$episodios.Range("F8").Value = "https://collaboration.merck.com/sites/onconceptos/_layouts/15/embed.aspx?UniqueId=3934dfb1-e42b-4782-8a02-093fbd8326a6&embed=%7B%22hvm%22%3Atrue%2C%22ust%22%3Atrue%7D&referrer=OneUpFileViewer&referrerScenario=EmbedDialog.Create"

# ---------------------------------------------------------------------------
# View state: selection moves to episodios!B8 while temporadas becomes the
# active (foreground) sheet with A3 selected.
# ---------------------------------------------------------------------------
$episodios.Activate()
$episodios.Range("A1").Select() | Out-Null
$episodios.Range("B8").Select() | Out-Null

$temporadas.Activate()
$temporadas.Range("A3").Select() | Out-Null
